$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.191.60"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").Value = "2.928.92"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.78"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.68"
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -0.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.00"
$ws.Range("E9").Value = "  +2.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").Value = "  -2.06%  "

$ws.Range("E11").Value = "  -1.02%  "

$ws.Range("E12").Value = "  -1.51%  "

$ws.Range("E13").Value = "  -0.98%  "

$ws.Range("E14").Value = "  +0.28%  "

$ws.Range("D15").Value = "3.412.86"
$ws.Range("E15").Value = "  -0.27%  "

$ws.Range("D16").Value = "61.193.97"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("E17").Value = "  -0.59%  "

$ws.Range("D18").Value = "2.924.91"
$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "432.67"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.48"
$ws.Range("E20").Value = "  -0.89%  "

$ws.Range("E21").Value = "  -1.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.12"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "82.00"
$ws.Range("E23").Value = "  +1.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.02"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("E25").Value = "  -1.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.80"
$ws.Range("E26").Value = "  -3.12%  "

$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("E28").Value = "  -5.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.60"
$ws.Range("E29").Value = "  -1.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").Value = "  -2.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.110"
$ws.Range("E31").Value = "  +1.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.81"
$ws.Range("E32").Value = "  +0.17%  "

$ws.Range("E33").Value = "  +0.24%  "

$ws.Range("D34").Value = "0.0₃0886"
$ws.Range("E34").Value = "  +1.52%  "

$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.65"
$ws.Range("E36").Value = "  -0.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.99"
$ws.Range("E37").Value = "  -4.94%  "

$ws.Range("E38").Value = "  -1.52%  "

$ws.Range("E39").Value = "  -1.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.62"
$ws.Range("E40").Value = "  -0.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.37"
$ws.Range("E41").Value = "  +4.88%  "

$ws.Range("E42").Value = "  -3.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0347"
$ws.Range("E43").Value = "  -1.24%  "

$ws.Range("D44").Value = "2.697.83"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "369.23"
$ws.Range("E45").Value = "  -3.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.50"
$ws.Range("E46").Value = "  +2.74%  "

$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.81"
$ws.Range("E48").Value = "  -2.23%  "

$ws.Range("E49").Value = "  -1.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.00"
$ws.Range("E50").Value = "  -1.62%  "

$ws.Range("E51").Value = "  -0.59%  "
